$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 11
$ws.Range("H11").Value = 1112
$ws.Range("I11").Value = 1112
$ws.Range("K11").Value = 1112
$ws.Range("M11").Value = -972
# Row 29
$ws.Range("H29").Value = 1544.6
$ws.Range("I29").Value = 1544.6
$ws.Range("K29").Value = 4633.799999999999
$ws.Range("M29").Value = -4352.799999999999
# Row 55
$ws.Range("H55").Value = 179.46153
$ws.Range("J55").Value = 267.69232
$ws.Range("L55").Value = 267.69232
$ws.Range("N55").Value = -695.69232
# Row 97
$ws.Range("H97").Value = 1646.5
$ws.Range("J97").Value = 2130.0588
$ws.Range("L97").Value = 6390.176399999999
$ws.Range("N97").Value = -7382.176399999999
# Row 107
$ws.Range("H107").Value = 30303684
$ws.Range("I107").Value = 66666812
$ws.Range("J107").Value = 1078
$ws.Range("K107").Value = 66666812
$ws.Range("L107").Value = 1078
$ws.Range("M107").Value = -66664892
$ws.Range("N107").Value = -4918
# Row 111
$ws.Range("H111").Value = 6176314.5
$ws.Range("I111").Value = 10103265
$ws.Range("J111").Value = 5392.857
$ws.Range("K111").Value = 30309795
$ws.Range("L111").Value = 16178.571
$ws.Range("M111").Value = -30306728
$ws.Range("N111").Value = -22312.571
# Row 118
$ws.Range("H118").Value = 100000570
$ws.Range("I118").Value = 111111620
$ws.Range("J118").Value = 1100
$ws.Range("K118").Value = 333334860
$ws.Range("L118").Value = 3300
$ws.Range("M118").Value = -333333203
$ws.Range("N118").Value = -6614
# Row 129
$ws.Range("H129").Value = 3992.3076
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 3992.3076
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 11976.9228
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -21976.9228
# Row 135
$ws.Range("H135").Value = 1568.525
$ws.Range("I135").Value = 770.44446
$ws.Range("J135").Value = 3226.077
$ws.Range("K135").Value = 6934.00014
$ws.Range("L135").Value = 29034.693
$ws.Range("M135").Value = -4399.00014
$ws.Range("N135").Value = -34104.693
# Row 137
$ws.Range("H137").Value = 42310.312
$ws.Range("I137").Value = 54379.207
$ws.Range("J137").Value = 5006.4546
$ws.Range("K137").Value = 163137.621
$ws.Range("L137").Value = 15019.3638
$ws.Range("M137").Value = -160587.621
$ws.Range("N137").Value = -20119.3638
# Row 138
$ws.Range("H138").Value = 3389.1025
$ws.Range("I138").Value = 2376.2942
$ws.Range("J138").Value = 3671.3606
$ws.Range("K138").Value = 7128.882599999999
$ws.Range("L138").Value = 11014.0818
$ws.Range("M138").Value = -1988.882599999999
$ws.Range("N138").Value = -21294.0818

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 5
$ws.Range("H5").Value = 204.08333
$ws.Range("I5").Value = 204.08333
$ws.Range("K5").Value = 204.08333
$ws.Range("M5").Value = -92.08332999999999
# Row 25
$ws.Range("H25").Value = 358
$ws.Range("I25").Value = 358
$ws.Range("K25").Value = 358
$ws.Range("M25").Value = 44
# Row 32
$ws.Range("H32").Value = 9718.644
$ws.Range("J32").Value = 28325
$ws.Range("L32").Value = 28325
$ws.Range("N32").Value = -28899
# Row 45
$ws.Range("H45").Value = 6805983.5
$ws.Range("I45").Value = 11905847
$ws.Range("K45").Value = 11905847
$ws.Range("M45").Value = -11905470
# Row 102
$ws.Range("H102").Value = 6947678
$ws.Range("I102").Value = 6947678
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 6947678
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -6946056
$ws.Range("N102").ClearContents()
# Row 110
$ws.Range("H110").Value = 1737884.9
$ws.Range("I110").Value = 1737884.9
$ws.Range("K110").Value = 1737884.9
$ws.Range("M110").Value = -1735839.9

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 204.08333
$ws.Range("I4").Value = 204.08333
$ws.Range("K4").Value = 204.08333
$ws.Range("M4").Value = -89.08332999999999
# Row 134
$ws.Range("H134").Value = 7162.8
$ws.Range("I134").Value = 2125.611
$ws.Range("K134").Value = 6376.833
$ws.Range("M134").Value = -3841.833

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 4
$ws.Range("H4").Value = 9000
$ws.Range("J4").Value = 9000
$ws.Range("L4").Value = 9000
$ws.Range("N4").Value = -9224
# Row 31
$ws.Range("H31").Value = 20878.59
$ws.Range("I31").Value = 3243.5173
$ws.Range("J31").Value = 39819.965
$ws.Range("K31").Value = 3243.5173
$ws.Range("L31").Value = 39819.965
$ws.Range("M31").Value = -2948.5173
$ws.Range("N31").Value = -40409.965
# Row 34
$ws.Range("H34").Value = 20878.59
$ws.Range("I34").Value = 3243.5173
$ws.Range("J34").Value = 39819.965
$ws.Range("K34").Value = 3243.5173
$ws.Range("L34").Value = 39819.965
$ws.Range("M34").Value = -3041.5173
$ws.Range("N34").Value = -40223.965
# Row 58
$ws.Range("H58").Value = 8624
$ws.Range("I58").Value = 10760.728
$ws.Range("J58").Value = 4706.6665
$ws.Range("K58").Value = 10760.728
$ws.Range("L58").Value = 4706.6665
$ws.Range("M58").Value = -10557.728
$ws.Range("N58").Value = -5112.6665
# Row 62
$ws.Range("H62").Value = 6046
$ws.Range("I62").Value = 5810.143
$ws.Range("K62").Value = 5810.143
$ws.Range("M62").Value = -5186.143
# Row 65
$ws.Range("H65").Value = 6046
$ws.Range("I65").Value = 5810.143
$ws.Range("K65").Value = 29050.715
$ws.Range("M65").Value = -25930.715
# Row 109
$ws.Range("H109").Value = 32129.5
$ws.Range("J109").Value = 34000
$ws.Range("L109").Value = 34000
$ws.Range("N109").Value = -36080
# Row 136
$ws.Range("H136").Value = 8624
$ws.Range("I136").Value = 10760.728
$ws.Range("J136").Value = 4706.6665
$ws.Range("K136").Value = 32282.184
$ws.Range("L136").Value = 14119.9995
$ws.Range("M136").Value = -29732.184
$ws.Range("N136").Value = -19219.9995
# Row 141
$ws.Range("H141").Value = 120304.86
$ws.Range("J141").Value = 120304.86
$ws.Range("L141").Value = 120304.86
$ws.Range("N141").Value = -130664.86

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 37
$ws.Range("H37").Value = 55250
$ws.Range("J37").Value = 55250
$ws.Range("L37").Value = 165750
$ws.Range("N37").Value = -165974
# Row 131
$ws.Range("H131").Value = 27785260
$ws.Range("J131").Value = 41679904
$ws.Range("L131").Value = 125039712
$ws.Range("N131").Value = -125049792
# Row 132
$ws.Range("H132").Value = 1753.2667
$ws.Range("I132").Value = 1650
$ws.Range("J132").Value = 1871.2858
$ws.Range("K132").Value = 14850
$ws.Range("L132").Value = 16841.5722
$ws.Range("M132").Value = -12320
$ws.Range("N132").Value = -21901.5722

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 22
$ws.Range("H22").Value = 15874.25
$ws.Range("J22").Value = 16999
$ws.Range("L22").Value = 16999
$ws.Range("N22").Value = -18057
# Row 52
$ws.Range("H52").Value = 30033
$ws.Range("J52").Value = 30033
$ws.Range("L52").Value = 30033
$ws.Range("N52").Value = -30551
# Row 132
$ws.Range("H132").Value = 2990.7542
$ws.Range("I132").Value = 2752.4255
$ws.Range("K132").Value = 8257.2765
$ws.Range("M132").Value = -5727.2765

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6290.7803
$ws.Range("I7").Value = 5759.3477
$ws.Range("J7").Value = 6969.8335
$ws.Range("K7").Value = 5759.3477
$ws.Range("L7").Value = 6969.8335
$ws.Range("M7").Value = -5647.3477
$ws.Range("N7").Value = -7193.8335
# Row 61
$ws.Range("H61").Value = 7938098.5
$ws.Range("I61").Value = 8548568
$ws.Range("K61").Value = 8548568
$ws.Range("M61").Value = -8548366
# Row 68
$ws.Range("H68").Value = 1000
$ws.Range("I68").Value = 1000
$ws.Range("K68").Value = 1000
$ws.Range("M68").Value = -251
# Row 71
$ws.Range("H71").Value = 1000
$ws.Range("I71").Value = 1000
$ws.Range("K71").Value = 5000
$ws.Range("M71").Value = -1256
# Row 81
$ws.Range("H81").Value = 47500
$ws.Range("J81").Value = 47500
$ws.Range("L81").Value = 47500
$ws.Range("N81").Value = -49496
# Row 84
$ws.Range("H84").Value = 47500
$ws.Range("J84").Value = 47500
$ws.Range("L84").Value = 142500
$ws.Range("N84").Value = -152484
# Row 113
$ws.Range("H113").Value = 7938098.5
$ws.Range("I113").Value = 8548568
$ws.Range("K113").Value = 8548568
$ws.Range("M113").Value = -8546398
# Row 122
$ws.Range("H122").Value = 5276.206
$ws.Range("I122").Value = 4253.8076
$ws.Range("J122").Value = 8599
$ws.Range("K122").Value = 12761.4228
$ws.Range("L122").Value = 25797
$ws.Range("M122").Value = -10311.4228
$ws.Range("N122").Value = -30697
# Row 126
$ws.Range("H126").Value = 6290.7803
$ws.Range("I126").Value = 5759.3477
$ws.Range("J126").Value = 6969.8335
$ws.Range("K126").Value = 17278.0431
$ws.Range("L126").Value = 20909.5005
$ws.Range("M126").Value = -14808.0431
$ws.Range("N126").Value = -25849.5005
# Row 133
$ws.Range("H133").Value = 198998
$ws.Range("J133").Value = 198998
$ws.Range("L133").Value = 198998
$ws.Range("N133").Value = -204058

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 100
$ws.Range("H100").Value = 4076.7334
$ws.Range("I100").Value = 4565.6924
$ws.Range("K100").Value = 9131.3848
$ws.Range("M100").Value = -8590.3848
# Row 126
$ws.Range("H126").Value = 2937.9
$ws.Range("I126").Value = 3399
$ws.Range("K126").Value = 10197
$ws.Range("M126").Value = -7727
# Row 136
$ws.Range("H136").Value = 2445.7847
$ws.Range("I136").Value = 2187.6345
$ws.Range("K136").Value = 6562.9035
$ws.Range("M136").Value = -4012.9035

Write-Host "Applied Hyperion_Profits scheduled-runner updates across all 8 sheets"
